# TC_69.xlsx: rename the worksheet and restyle the column-B values.
#
#   - Worksheet "My Series" -> "Data"
#   - Column B (the observation values) drops its bold / red font and
#     goes back to the workbook's default (theme) font
#   - Column B's custom number format becomes "###0.0000" (was "0.0000")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet.
$ws.Name = "Data"

# 2. Re-style the data column (B1:B10).
$rng = $ws.Range("B1:B10")

# Drop the bold weight and the custom red color, returning to the
# default theme-colored, regular-weight font used elsewhere in the sheet.
$rng.Font.FontStyle = "Regular"
$rng.Font.ThemeColor = 1

# Widen the custom numeric format.
$rng.NumberFormat = "###0.0000"
